$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking columns (G:K) to be stored as text, matching the
# source data (which keeps values like "140.00" / leading context as text).
$ws.Range("G2:K9").NumberFormat = "@"

$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"
$ws.Range("A2").Value = " Dubai (DSC)"
$ws.Range("B2").Value = " October 24 2020"
$ws.Range("C2").Value = "Kings XI won by 12 runs"
$ws.Range("D2").Value = "Sunrisers Hyderabad"
$ws.Range("E2").Value = "Kings XI Punjab"
$ws.Range("F2").Value = "Abdul Samad "
$ws.Range("G2").Value = "7"
$ws.Range("H2").Value = "5"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "140.00"
$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " October 18 2020"
$ws.Range("C3").Value = "Match tied (KKR won the one-over eliminator)"
$ws.Range("D3").Value = "Sunrisers Hyderabad"
$ws.Range("E3").Value = "Kolkata Knight Riders"
$ws.Range("F3").Value = "Abdul Samad "
$ws.Range("G3").Value = "23"
$ws.Range("H3").Value = "15"
$ws.Range("I3").Value = "2"
$ws.Range("J3").Value = "1"
$ws.Range("K3").Value = "153.33"
$ws.Range("A4").Value = " Sharjah"
$ws.Range("B4").Value = " October 31 2020"
$ws.Range("C4").Value = "Sunrisers won by 5 wickets (with 35 balls remaining)"
$ws.Range("D4").Value = "Sunrisers Hyderabad"
$ws.Range("E4").Value = "Royal Challengers Bangalore"
$ws.Range("F4").Value = "Abdul Samad "
$ws.Range("G4").Value = "0"
$ws.Range("H4").Value = "0"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "-"
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " November 08 2020"
$ws.Range("C5").Value = "Capitals won by 17 runs"
$ws.Range("D5").Value = "Sunrisers Hyderabad"
$ws.Range("E5").Value = "Delhi Capitals"
$ws.Range("F5").Value = "Abdul Samad "
$ws.Range("G5").Value = "33"
$ws.Range("H5").Value = "16"
$ws.Range("I5").Value = "2"
$ws.Range("J5").Value = "2"
$ws.Range("K5").Value = "206.25"
$ws.Range("A6").Value = " Dubai (DSC)"
$ws.Range("B6").Value = " October 08 2020"
$ws.Range("C6").Value = "Sunrisers won by 69 runs"
$ws.Range("D6").Value = "Sunrisers Hyderabad"
$ws.Range("E6").Value = "Kings XI Punjab"
$ws.Range("F6").Value = "Abdul Samad "
$ws.Range("G6").Value = "8"
$ws.Range("H6").Value = "7"
$ws.Range("I6").Value = "1"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "114.28"
$ws.Range("A7").Value = " Dubai (DSC)"
$ws.Range("B7").Value = " October 02 2020"
$ws.Range("C7").Value = "Sunrisers won by 7 runs"
$ws.Range("D7").Value = "Sunrisers Hyderabad"
$ws.Range("E7").Value = "Chennai Super Kings"
$ws.Range("F7").Value = "Abdul Samad "
$ws.Range("G7").Value = "8"
$ws.Range("H7").Value = "6"
$ws.Range("I7").Value = "0"
$ws.Range("J7").Value = "0"
$ws.Range("K7").Value = "133.33"
$ws.Range("A8").Value = " Sharjah"
$ws.Range("B8").Value = " October 04 2020"
$ws.Range("C8").Value = "Mumbai won by 34 runs"
$ws.Range("D8").Value = "Sunrisers Hyderabad"
$ws.Range("E8").Value = "Mumbai Indians"
$ws.Range("F8").Value = "Abdul Samad "
$ws.Range("G8").Value = "20"
$ws.Range("H8").Value = "9"
$ws.Range("I8").Value = "1"
$ws.Range("J8").Value = "2"
$ws.Range("K8").Value = "222.22"
$ws.Range("A9").Value = " Abu Dhabi"
$ws.Range("B9").Value = " September 29 2020"
$ws.Range("C9").Value = "Sunrisers won by 15 runs"
$ws.Range("D9").Value = "Sunrisers Hyderabad"
$ws.Range("E9").Value = "Delhi Capitals"
$ws.Range("F9").Value = "Abdul Samad "
$ws.Range("G9").Value = "12"
$ws.Range("H9").Value = "7"
$ws.Range("I9").Value = "1"
$ws.Range("J9").Value = "1"
$ws.Range("K9").Value = "171.42"
